$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with Southbank / Rockpool exposure site (old period)
$ws.Range("A2").Value = "Southbank"
$ws.Range("B2").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C2").Value = "23/12/2020 1:00pm-1:30pm"
$ws.Range("D2").Value = "Case attended restaurant"
$ws.Range("E2").Value = "old"

# Add row 3 with Southbank / Rockpool exposure site (new period)
$ws.Range("A3").Value = "Southbank"
$ws.Range("B3").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C3").Value = "23/12/2020 8:00pm-10:00pm"
$ws.Range("D3").Value = "Case attended restaurant"
$ws.Range("E3").Value = "new"
